$d = $word.ActiveDocument

# --- Step 1: remove the "Meta description" paragraph near the top of the
#     document (2nd paragraph): empty run + bold "Meta description" run +
#     plain ": Discover the ancient secrets ..." run. The new paragraph we
#     add later reuses this exact same run-shape, so grab its formatted
#     contents before deleting it.
$metaPara = $d.Paragraphs.Item(2)
$metaFormattedText = $metaPara.Range.FormattedText
$metaPara.Range.Delete()

# --- Step 2: insert a new paragraph, with the same 3-run shape as the
#     "Meta description" paragraph (empty run / bold run / plain run),
#     directly before the final "Prompt for DALLE ..." paragraph. We build
#     it right after the paragraph before last ("No scatter or bonus
#     symbols"), which has plain (non-bold/non-italic) trailing formatting,
#     so the freshly created paragraph/run picks up clean formatting.
$n = $d.Paragraphs.Count
$anchor = $d.Paragraphs.Item($n - 1)
$anchor.Range.InsertParagraphAfter()

$newPara = $d.Paragraphs.Item($n)
$newPara.Range.ParagraphFormat.Style = "Normal"

$newPara2 = $d.Paragraphs.Item($n)
$newPara2.Range.FormattedText = $metaFormattedText

# Collapse the copied paragraph's run content (bold "Meta description" +
# plain ": Discover ...") down to a single run with the new heading text,
# leaving the genuinely empty leading run intact.
$newPara3 = $d.Paragraphs.Item($n)
$start = $newPara3.Range.Start
$end = $newPara3.Range.End
$textRange = $d.Range($start, $end - 1)
$textRange.Text = "Play Black Mummy Free - Slot Game Review | Tom Horn Gaming"

# --- Step 3: update the (now last) "Prompt for DALLE ..." paragraph's
#     text, keeping its existing empty leading run + italic run formatting.
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$lastStart = $lastPara.Range.Start
$lastEnd = $lastPara.Range.End
$lastTextRange = $d.Range($lastStart, $lastEnd - 1)
$lastTextRange.Text = "Discover the ancient secrets of Black Mummy and play for free. Read our review of Tom Horn Gaming's online slot machine to uncover its gameplay mechanics and features."
